# [Improvement] On terminology : room -> bed
$wb = $excel.ActiveWorkbook

$babies = $wb.Worksheets.Item("babies")
$rooms  = $wb.Worksheets.Item("rooms")

# Rename "rooms" sheet to "beds" and update the room->bed terminology
# used in the header row of that sheet.
$rooms.Range("A1").Value = "all_beds"
$rooms.Range("B1").Value = "new_beds"
$rooms.Range("C1").Value = "old_beds"
$rooms.Range("E1").Value = "new_beds_service"
$rooms.Range("F1").Value = "old_beds_service"
$rooms.Range("G1").Value = "beds_capacities"
$rooms.Name = "beds"

# Update selections and the active sheet/tab.
$babies.Activate() | Out-Null
$babies.Range("E13").Select() | Out-Null

$rooms.Activate() | Out-Null
$rooms.Range("I2").Select() | Out-Null
